# Applies the OOP-in-BF content rewrite: reworks the intro paragraphs, replaces the
# "It is very limited in its scope." ending with a fuller description of BF's command
# set, and appends a new "core features" list (Inheritance / Polymorphism /
# Encapsulation / Abstraction) plus trailing blank paragraphs.
$d = $word.ActiveDocument

# --- Paragraph 3: "Let us start by defining what an object is..." ---
# split the sentence around the old "object is" clause into 3 runs, swapping in the
# new "makes a programming language object orientated" clause.
$d.Paragraphs(3).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Let us start by defining what an </w:t></w:r><w:r><w:t>makes a programming language object orientated</w:t></w:r><w:r><w:t xml:space="preserve"> and how it could potentially be applied to BF. </w:t></w:r></w:p></w:body></w:document>') | Out-Null

# --- Paragraph 5: "consists of state and related behaviour..." ---
# replaced outright with the old paragraph-3 tail sentence ("An object is an
# abstract data type...").
$d.Paragraphs(5).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">An object is an abstract data type created by a developer. It can include multiple properties (state) and methods. In most programming languages, objects are defined as classes. </w:t></w:r></w:p></w:body></w:document>') | Out-Null

# --- Paragraph 6: "BF is a programming language built on eight simple commands..." ---
# drop the "It is very limited in its scope." tail and replace it with a longer,
# multi-run description of what the commands do.
$d.Paragraphs(6).Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">BF is a programming language built on eight simple commands and an instruction pointer. </w:t></w:r><w:r><w:t xml:space="preserve">These </w:t></w:r><w:r><w:t>commands revolve around manipulating the stack</w:t></w:r><w:r><w:t xml:space="preserve"> and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">printing </w:t></w:r><w:r><w:t>the ASCII representation of the value stored in the current pointer location</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document>') | Out-Null

# --- New paragraphs appended after paragraph 6 ---
function Add-ParagraphWithXml([string]$xml) {
    $count = $d.Paragraphs.Count
    $tail = $d.Paragraphs($count).Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter() | Out-Null
    $d.Paragraphs($count + 1).Range.InsertXML($xml) | Out-Null
}

Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The core </w:t></w:r><w:r><w:t xml:space="preserve">features that need to be applied to BF to </w:t></w:r></w:p></w:body></w:document>'  # "The core features that need to be applied to BF to "
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Inheritance</w:t></w:r></w:p></w:body></w:document>'  # "Inheritance"
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Polymorphism</w:t></w:r></w:p></w:body></w:document>'  # "Polymorphism"
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Encapsulation</w:t></w:r></w:p></w:body></w:document>'  # "Encapsulation"
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Abstraction</w:t></w:r></w:p></w:body></w:document>'  # "Abstraction"
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p></w:p></w:body></w:document>'  # trailing blank paragraph
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p></w:p></w:body></w:document>'  # trailing blank paragraph
Add-ParagraphWithXml '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p></w:p></w:body></w:document>'  # trailing blank paragraph

